# Apply the LOT2017.xlsx content update (course info reshuffle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOT2017"
$ws.Range("C2").Value = "LOT2017"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Enzimologia"
$ws.Range("C3").Value = " Enzimologia"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Enzymology"
$ws.Range("C4").Value = "Enzymology"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "3"
$ws.Range("C5").Value = "3"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "45 h"
$ws.Range("C7").Value = "45 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2018"
$ws.Range("C8").Value = "01/01/2018"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EB-7"
$ws.Range("C9").Value = "EB-7"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "4873328 - Fernando Segato"
$ws.Range("C10").Value = "4873328 - Fernando Segato"

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "The discipline deals with how the enzymes act, how the strategies of purification are defined and what are the main technological applications of the enzymes. Within the purification processes, the focus involves the definition of appropriate strategies for the purification in sequential stages, the control methods of each stage, besides the methods of monitoring of the enzymatic activity. The application of enzymes in industrial processes is also commented."
$ws.Range("C14").Value = "The discipline deals with how the enzymes act, how the strategies of purification are defined and what are the main technological applications of the enzymes. Within the purification processes, the focus involves the definition of appropriate strategies for the purification in sequential stages, the control methods of each stage, besides the methods of monitoring of the enzymatic activity. The application of enzymes in industrial processes is also commented."

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1.Cellular origin of enzymes: enzymes origin, differentiation between intra and extracellular enzymes, physiological importance and introduction to the world market of enzymes.2.Structure versus proprieties and mechanisms of enzymes action; three-dimensional structure and its determination, importance of tertiary structure on the catalytic activity, catalytic action of protease, glycosidase and oxy-reductase.3.Operational control in the enzymes purification: methods of enzymes extraction, methods of preliminary purification, methods of separation based on charge, size and affinity. Definition of strategies of enzymes purification.4.Methods of determination of enzymatic activity: definition of activity in international units, meaning of enzymatic activity, forms of quantifying and expressing enzymatic activity. Requirements of an experimental method used in the determination of enzymatic activity.5.Enzymatic kinetics: graphic and numeric methods for determination of initial reaction rate, experimental conditions required to determine the initial reaction rate, calculation of enzymatic activity.6.Immobilized enzymes: forms of immobilization and application of immobilized systems.7.Application of enzymes in industry: use of enzymes in detergents, starch processing, food industry, textile industry, drug synthesis and in pulp and paper industry."
$ws.Range("C16").Value = "1.Cellular origin of enzymes: enzymes origin, differentiation between intra and extracellular enzymes, physiological importance and introduction to the world market of enzymes.2.Structure versus proprieties and mechanisms of enzymes action; three-dimensional structure and its determination, importance of tertiary structure on the catalytic activity, catalytic action of protease, glycosidase and oxy-reductase.3.Operational control in the enzymes purification: methods of enzymes extraction, methods of preliminary purification, methods of separation based on charge, size and affinity. Definition of strategies of enzymes purification.4.Methods of determination of enzymatic activity: definition of activity in international units, meaning of enzymatic activity, forms of quantifying and expressing enzymatic activity. Requirements of an experimental method used in the determination of enzymatic activity.5.Enzymatic kinetics: graphic and numeric methods for determination of initial reaction rate, experimental conditions required to determine the initial reaction rate, calculation of enzymatic activity.6.Immobilized enzymes: forms of immobilization and application of immobilized systems.7.Application of enzymes in industry: use of enzymes in detergents, starch processing, food industry, textile industry, drug synthesis and in pulp and paper industry."

$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).EntireRow.AutoFit()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "4873328 - Fernando Segato"
$ws.Range("C18").Value = "4873328 - Fernando Segato"
$ws.Rows.Item(18).RowHeight = 60

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas (P1 e P2)."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas (P1 e P2)."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Rows.Item(21).RowHeight = 120

$ws.Range("A22").Value = "Requisitos:"
$ws.Rows.Item(22).EntireRow.AutoFit()

$ws.Range("B23").Value = "LOT2007 -  Bioquímica I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOT2007 -  Bioquímica I  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("B24").Value = "LOT2040 -  Engenharia Genética  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOT2040 -  Engenharia Genética  (Requisito fraco)`n"

$ws.Range("B25").Value = "LOT2053 -  Microbiologia  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOT2053 -  Microbiologia  (Requisito fraco)`n"

# Row 26 (old "LOT2053..." row) is removed entirely; everything above already holds its place
$ws.Rows.Item(26).Delete()

Write-Output "applied LOT2017 content update"